$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 31 continues the existing schedule pattern:
#   A column = date 9 rows back (reuses the existing "10/24" shared string)
#   B column = new period date label
#   C column = new period description
$ws.Range("A31").Value = "10/24"
$ws.Range("B31").Value = "12/19"
$ws.Range("C31").Value = "第81期 混合紙飛機 輪次獲得6鈴鐺 750能量 250貓木  12坐騎碎片  三層鈴鐺5個 四層 20自選秘寶碎或15麵粉  五層 神話坐騎5點可換兩次"

# Match the existing style of column A/B cells (text number format)
$ws.Range("A31:B31").NumberFormat = "@"

$ws.Range("C31").Select()
